$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold numeric-looking strings (e.g. "68.077.62", "  +0.33%  ")
# that must stay literal text (matching the source inlineStr cells), so
# force text format before assigning to avoid Excel auto-converting them
# to numbers/percentages and losing precision or formatting.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.077.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.684.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.04"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.99"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.709"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +12.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000277"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.94%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.23"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.50%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.292.25"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.40%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.698.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.91%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.61"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.009.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "402.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.48"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.71"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.56"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.74%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.99"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.64"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.74%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.99"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.52%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.63"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.30%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.42"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.10%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "45.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.58%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "67.21"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.118"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "619.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.66%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.28%  "

# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.401"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0792"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -9.64%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.29%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.97%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0431"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.72%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.862.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.137"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.43%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.11%  "

# Row 48
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.05"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.85%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.65"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.07%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.03"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.46%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -11.17%  "
